# Jeannine's log - add Thursday 9/8/2016 entries (CLH rooms G & I, neck mic/PC demos)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 137 is a day-separator row (same look as the other "THURSDAY" banners,
# e.g. row 5). Copy that formatting onto the new row, then set its label.
$ws.Range("A5:F5").Copy()
$ws.Range("A137:F137").PasteSpecial(-4122)
$ws.Range("B137").Value = "THURSDAY"

# Row 138: Demo / 09-08-2016 / 1630 / CLH / Room G / comment
$ws.Range("A138").Value = "Demo"
$ws.Range("B138").Value = 42621
$ws.Range("C138").Value = "1630"
$ws.Range("D138").Value = "CLH"
$ws.Range("E138").Value = "G"
$ws.Range("F138").Value = "Demo neck mic and PC"

# Row 139: Demo / 09-08-2016 / 1850 / CLH / Room I / comment (trailing space kept)
$ws.Range("A139").Value = "Demo"
$ws.Range("B139").Value = 42621
$ws.Range("C139").Value = "1850"
$ws.Range("D139").Value = "CLH"
$ws.Range("E139").Value = "I"
$ws.Range("F139").Value = "Demo neck mic and PC "

# Keep the view positioned near the bottom of the log, same cell selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 115
$win.ScrollColumn = 1
$ws.Range("F139").Select() | Out-Null
